$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13: "Formula Test" scenario example row
$ws.Range("B13").Value = "Formula Test"
$ws.Range("C13").Value = "FormulaTest.xlsx"
$ws.Range("D13").Formula = "=SUBSTITUTE(C13, "".xlsx"", "".feature"")"
$ws.Range("E13").Formula = "=SUBSTITUTE(C13, "".xlsx"", "".exp"")"

# New cell D6: "except empty lines"
$ws.Range("D6").Value = "except empty lines"

# C5 loses its bold/font style override - reset to default style
$ws.Range("C5").Style = "Normal"

# Move the active selection to D18 (matches diff's new selection)
$ws.Range("D18").Select() | Out-Null

# Enable iterative calculation (workbook.xml calcPr iterateCount/iterateDelta)
$excel.Iteration = $true
$excel.MaxIterations = 250
$excel.MaxChange = 0.00001
